# Update countries & provincias Spain
# Applies:
#  1) Six pairs of adjacent countries in the ranking swap places (both the
#     country name and -- where applicable -- their stats), reflecting a
#     re-sort of the "Pais" table by total cases.
#  2) A refresh of the COVID-19 case numbers for a number of countries.
#  3) The "last updated" timestamp string in cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap country-name pairs (rows that traded ranking positions) -------
$ws.Range("A103").Value = "Costa Rica"
$ws.Range("A104").Value = "Maldivas"

$ws.Range("A142").Value = "Ruanda"
$ws.Range("A143").Value = "Crucero"

$ws.Range("A168").Value = "Angola"
$ws.Range("A169").Value = "Guyana"

$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A207").Value = "Islas Malvinas"
$ws.Range("A208").Value = "Groenlandia"

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- 2) Refresh updated statistics ------------------------------------------
$ws.Range("B4").Value = 2347701
$ws.Range("C4").Value = 17123
$ws.Range("E4").Value = 1250489
$ws.Range("G4").Value = 183
$ws.Range("H4").Value = 122163
$ws.Range("B7").Value = 426910
$ws.Range("C7").Value = 15183
$ws.Range("D7").Value = 237252
$ws.Range("E7").Value = 175955
$ws.Range("G7").Value = 426
$ws.Range("H7").Value = 13703
$ws.Range("B14").Value = 191346
$ws.Range("C14").Value = 130
$ws.Range("E14").Value = 7484
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 8962
$ws.Range("B18").Value = 160377
$ws.Range("C18").Value = 284
$ws.Range("E18").Value = 56425
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 29640
$ws.Range("B22").Value = 97302
$ws.Range("C22").Value = 4621
$ws.Range("D22").Value = 51608
$ws.Range("E22").Value = 43764
$ws.Range("G22").Value = 53
$ws.Range("H22").Value = 1930
$ws.Range("B29").Value = 55233
$ws.Range("C29").Value = 1475
$ws.Range("D29").Value = 14736
$ws.Range("E29").Value = 38304
$ws.Range("G29").Value = 87
$ws.Range("H29").Value = 2193
$ws.Range("E50").Value = 5479
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 62
$ws.Range("B76").Value = 6315
$ws.Range("C76").Value = 162
$ws.Range("D76").Value = 4377
$ws.Range("E76").Value = 1919
$ws.Range("B96").Value = 2984
$ws.Range("C96").Value = 171
$ws.Range("D96").Value = 811
$ws.Range("E96").Value = 2062
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 111
$ws.Range("B103").Value = 2213
$ws.Range("C103").Value = 86
$ws.Range("D103").Value = 1032
$ws.Range("E103").Value = 1169
$ws.Range("H103").Value = 12
$ws.Range("B104").Value = 2203
$ws.Range("C104").Value = 16
$ws.Range("D104").Value = 1803
$ws.Range("E104").Value = 392
$ws.Range("H104").Value = 8
$ws.Range("B142").Value = 728
$ws.Range("C142").Value = 26
$ws.Range("D142").Value = 359
$ws.Range("E142").Value = 367
$ws.Range("H142").Value = 2
$ws.Range("B143").Value = 712
$ws.Range("D143").Value = 651
$ws.Range("E143").Value = 48
$ws.Range("H143").Value = 13
$ws.Range("B148").Value = 635
$ws.Range("C148").Value = 8
$ws.Range("E148").Value = 345
$ws.Range("B150").Value = 569
$ws.Range("C150").Value = 8
$ws.Range("D150").Value = 375
$ws.Range("E150").Value = 181
$ws.Range("B154").Value = 489
$ws.Range("C154").Value = 10
$ws.Range("D154").Value = 64
$ws.Range("E154").Value = 419
$ws.Range("G154").Value = 2
$ws.Range("H154").Value = 6
$ws.Range("C168").Value = 7
$ws.Range("D168").Value = 77
$ws.Range("E168").Value = 97
$ws.Range("H168").Value = 9
$ws.Range("B169").Value = 183
$ws.Range("D169").Value = 102
$ws.Range("E169").Value = 69
$ws.Range("H169").Value = 12
$ws.Range("B180").Value = 100
$ws.Range("C180").Value = 1
$ws.Range("E180").Value = 2
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

# --- 3) Update the "last updated" timestamp ---------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 21:38"
